$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "sample task"
$ws.Range("C10").Value = "Sample Images/download_figma_and_install.PNG"

$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "sample task1"
$ws.Range("C11").Value = "Sample Images/download_figma_and_install.PNG"
